{"js": "// Apply the benchmark-stats edit: rewrite specific cells of the (single\n// column) results table. Row indices below are 0-based, matching the\n// table's row order top-to-bottom.\nconst body = context.document.body;\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// row index (0-based) -> new text\nconst edits = {\n  0: \"0M\",       // was \"100\"\n  1: \"0M\",       // was \"0.05\"\n  2: \"0M\",       // was \"3864\"\n  3: \"863\",      // was \"783\"\n  8: \"0.00008\",  // was \"0.00005\"\n  11: \"0.04946\", // was \"0.04471\"\n  43: \"100\",     // was the multi-run/tabbed \"76\\t0.00003\\t...\\t100.0\" row\n  44: \"0.05\",    // was the multi-run/tabbed \"1\\t0.00008\\t...\\t100.0\" row\n  45: \"3864\",    // was the multi-run/tabbed \"3\\t0.00003\\t...\\t100.0\" row\n};\n\nfor (const [row, text] of Object.entries(edits)) {\n  table.getCell(Number(row), 0).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Apply the benchmark-stats edit: rewrite specific cells of the (single\n# column) results table. Word COM Table.Cell(row, col) is 1-based.\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$tbl.Cell(1, 1).Range.Text = \"0M\"        # was \"100\"\n$tbl.Cell(2, 1).Range.Text = \"0M\"        # was \"0.05\"\n$tbl.Cell(3, 1).Range.Text = \"0M\"        # was \"3864\"\n$tbl.Cell(4, 1).Range.Text = \"863\"       # was \"783\"\n$tbl.Cell(9, 1).Range.Text = \"0.00008\"   # was \"0.00005\"\n$tbl.Cell(12, 1).Range.Text = \"0.04946\"  # was \"0.04471\"\n$tbl.Cell(44, 1).Range.Text = \"100\"      # was the multi-run/tabbed \"76`t0.00003`t...`t100.0\" row\n$tbl.Cell(45, 1).Range.Text = \"0.05\"     # was the multi-run/tabbed \"1`t0.00008`t...`t100.0\" row\n$tbl.Cell(46, 1).Range.Text = \"3864\"     # was the multi-run/tabbed \"3`t0.00003`t...`t100.0\" row\n"}
